# Update workbook to reflect data refresh through 2022-08-02 (commit: "Add data for 2022-08-10")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab: "Through 2022-08-01" -> "Through 2022-08-02"
$ws.Name = "Through 2022-08-02"

# Update the header label in I1 (shared string "2022 (through 08-01)" -> "2022 (through 08-02)")
$ws.Range("I1").Value = "2022 (through 08-02)"

# Updated monthly figures for the "2022 (through ...)" column (I)
$ws.Range("I8").Value = 166   # July
$ws.Range("I9").Value = 6     # August

# Updated 2021 (column H) December figure
$ws.Range("H13").Value = 206  # December 2021

# Updated Total row (row 14)
$ws.Range("H14").Value = 1849 # 2021 total
$ws.Range("I14").Value = 977  # 2022 (through ...) total
